$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the volatile TIME(RANDBETWEEN(...)) formulas in column E (E2:E16)
# with their frozen, static computed values.
$ws.Range("E2").Value  = 0.62212962962962959
$ws.Range("E3").Value  = 0.43752314814814813
$ws.Range("E4").Value  = 0.96317129629629628
$ws.Range("E5").Value  = 0.65208333333333335
$ws.Range("E6").Value  = 0.067604166666666674
$ws.Range("E7").Value  = 0.71983796296296299
$ws.Range("E8").Value  = 0.43186342592592591
$ws.Range("E9").Value  = 0.10357638888888888
$ws.Range("E10").Value = 0.43119212962962961
$ws.Range("E11").Value = 0.09087962962962963
$ws.Range("E12").Value = 0.32666666666666666
$ws.Range("E13").Value = 0.24042824074074073
$ws.Range("E14").Value = 0.77543981481481483
$ws.Range("E15").Value = 0.64909722222222221
$ws.Range("E16").Value = 0.35109953703703706

# Update the selection shown in the saved worksheet view
$ws.Range("E2:E16").Select()
